$d = $word.ActiveDocument

# Locate the paragraph that ends the "Semaine 10" list
# ("Séparation du framework, et de l'application (js, imports, css)")
# and append a new list item with the same style/numbering right after it.
$target = $null
$targetIndex = -1
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    if ($para.Range.Text -like "*Séparation du framework, et de l'application (js, imports, css)*") {
        $target = $para
        $targetIndex = $i
    }
}

# InsertParagraphAfter clones the source paragraph's style/numbering
# (pStyle "Paragraphedeliste", ilvl 0, numId 1) onto the new paragraph.
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Réparation du champs autocomplete (probleme de nom de variable)"
